$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "51.724.53"
$ws.Range("E2").Value = "  +4.45%  "
$ws.Range("D3").Value = "2.767.01"
$ws.Range("E3").Value = "  +5.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.538"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +6.10%  "
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("D15").Value = "3.200.88"
$ws.Range("E15").Value = "  +5.25%  "
$ws.Range("D16").Value = "2.775.39"
$ws.Range("E16").Value = "  +5.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").Value = "51.641.76"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  +5.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0821"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +4.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0356"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("E44").Value = "  +7.50%  "
$ws.Range("E45").Value = "  +16.58%  "
$ws.Range("D46").Value = "2.090.17"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("E49").Value = "  +6.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
